$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 82

# Columns A-D hold text-like values (dates/times/weekday/week stored as text
# in this sheet, not as real Excel dates/numbers). Force text formatting so
# Excel doesn't auto-convert "2025-02-21" to a date serial or "07" to 7.
$ws.Range("A82:D82").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-02-21"
$ws.Cells.Item($row, 2).Value = "23:04:48"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "07"
$ws.Cells.Item($row, 5).Value = 130460
$ws.Cells.Item($row, 6).Value = 141333
$ws.Cells.Item($row, 7).Value = 172250
$ws.Cells.Item($row, 8).Value = 157241
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 146564
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193408
$ws.Cells.Item($row, 14).Value = 115377
$ws.Cells.Item($row, 15).Value = 46218
$ws.Cells.Item($row, 16).Value = 29293
$ws.Cells.Item($row, 17).Value = 68223
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48720
$ws.Cells.Item($row, 20).Value = -1
